$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: convert to "group end" style (add A12 border cell, switch styles 4/5 -> 6/7) ---
$ws.Range("A6:E6").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)

# --- Set B (row-number) values for the new rows first (plain numbers, no shared-string impact) ---
$ws.Range("B13").Value = 489
$ws.Range("B14").Value = 492
$ws.Range("B15").Value = 443
$ws.Range("B16").Value = 446
$ws.Range("B17").Value = 449
$ws.Range("B18").Value = 458
$ws.Range("B19").Value = 467
$ws.Range("B20").Value = 470
$ws.Range("B21").Value = 418
$ws.Range("B22").Value = 421
$ws.Range("B23").Value = 424
$ws.Range("B24").Value = 393
$ws.Range("B25").Value = 396

# --- Set C/D/E text values in the exact order the original authoring tool created them, ---
# --- so newly-appended sharedStrings entries land at the same indices as the target file. ---
$ws.Range("C13").Value = ' Meh-heh-heh. That aroma\nwafting down from upstairs…'
$ws.Range("C14").Value = ' I don\''t mind it at all. Nope, it\''s\nquite all right with me. Meh-heh-heh.'
$ws.Range("D13").Value = ' Хе-хе-хе. Этот аромат,\nдоносящийся сверху...'
$ws.Range("D14").Value = ' Я не возражаю его нюхать.\nНет, я даже не против него. Хе-хе-хе.'
$ws.Range("E13").Value = ' Öå-öå-öå. Üóïó àñïíàó,\näïîïòÿþéêòÿ òâåñöô...'
$ws.Range("E14").Value = ' Ÿ îå âïèñàçàý åãï îýöàóû.\nÎåó, ÿ äàçå îå ðñïóéâ îåãï. Öå-öå-öå.'
$ws.Range("C15").Value = ' Meh-heh-heh. I guess I get to\nenjoy this fine, ripe aroma for a while yet.'
$ws.Range("C16").Value = ' It\''s another thing that makes me\nhappy. Meh-heh-heh.'
$ws.Range("C17").Value = ' Oh, by the way…'
$ws.Range("C18").Value = ' I\''m almost done repairing my\nSwap Cauldron here…'
$ws.Range("C19").Value = ' When I\''m done fixing it, I\''ll\nbe sure to let you know.'
$ws.Range("C20").Value = ' You\''ll finally learn what I\''m up\nto here. Meh-heh-heh.'
$ws.Range("D15").Value = ' Хе-хе-хе. Кажется, я уже очень\nдавно наслаждаюсь этим чудесным, сочным\nароматом.'
$ws.Range("D16").Value = ' Он одна из немногих вещей, что\nделают меня счастливым. Хе-хе-хе.'
$ws.Range("D17").Value = ' О, кстати...'
$ws.Range("D18").Value = ' Я почти починил свой Обменный\nКотёл...'
$ws.Range("D19").Value = ' Когда я закончу его чинить, я\nдам вам об этом знать.'
$ws.Range("D20").Value = ' Наконец-то вы узнаете зачем я\nздесь. Хе-хе-хе.'
$ws.Range("E15").Value = ' Öå-öå-öå. Ëàçåóòÿ, ÿ ôçå ïœåîû\näàâîï îàòìàçäàýòû üóéí œôäåòîúí, òïœîúí\nàñïíàóïí.'
$ws.Range("E16").Value = ' Ïî ïäîà éè îåíîïãéö âåþåê, œóï\näåìàýó íåîÿ òœàòóìéâúí. Öå-öå-öå.'
$ws.Range("E17").Value = ' Ï, ëòóàóé...'
$ws.Range("E18").Value = ' Ÿ ðïœóé ðïœéîéì òâïê Ïáíåîîúê\nËïóæì...'
$ws.Range("E19").Value = ' Ëïãäà ÿ èàëïîœô åãï œéîéóû, ÿ\näàí âàí ïá üóïí èîàóû.'
$ws.Range("E20").Value = ' Îàëïîåø-óï âú ôèîàåóå èàœåí ÿ\nèäåòû. Öå-öå-öå.'
$ws.Range("C21").Value = ' I\''m almost done repairing my\nSwap Cauldron here... Almost, but not yet.'
$ws.Range("C22").Value = ' I\''ll tell you when I\''m done. It\''s\nsomething you should look forward to.\nMeh-heh-heh.'
$ws.Range("D21").Value = ' Я почти закончил чинить свой\nОбменный Котёл... Почти, но это ещё не\nвсё.'
$ws.Range("D22").Value = ' Когда я закончу, я вам сообщу.\nЖдите этого. Хе-хе-хе.'
$ws.Range("E21").Value = ' Ÿ ðïœóé èàëïîœéì œéîéóû òâïê\nÏáíåîîúê Ëïóæì... Ðïœóé, îï üóï åþæ îå\nâòæ.'
$ws.Range("E22").Value = ' Ëïãäà ÿ èàëïîœô, ÿ âàí òïïáþô.\nÇäéóå üóïãï. Öå-öå-öå.'
$ws.Range("C24").Value = ' We\''re closing up the guild while\nwe\''re on the expedition. I have to leave my\nSwap Cauldron behind. Meh-heh-heh.'
$ws.Range("C25").Value = ' I can\''t do this while on the\nexpedition, so if you have anything to swap,\ndo it now. Meh-heh-heh.'
$ws.Range("D24").Value = ' Во время экспедиции гильдия\nбудет закрыта. Мне придётся оставить\nздесь свой Обменный Котёл. Хе-хе-хе.'
$ws.Range("D25").Value = ' В экспедиции я ничего не смогу\nобменять, поэтому если вам есть что\nменять, делайте это сейчас. Хе-хе-хе.'
$ws.Range("E24").Value = ' Âï âñåíÿ üëòðåäéøéé ãéìûäéÿ\náôäåó èàëñúóà. Íîå ðñéäæóòÿ ïòóàâéóû\nèäåòû òâïê Ïáíåîîúê Ëïóæì. Öå-öå-öå.'
$ws.Range("E25").Value = ' Â üëòðåäéøéé ÿ îéœåãï îå òíïãô\nïáíåîÿóû, ðïüóïíô åòìé âàí åòóû œóï\níåîÿóû, äåìàêóå üóï òåêœàò. Öå-öå-öå.'

# --- Row 23 reuses the same strings as row 20 (dedup; no new sharedStrings entries) ---
$ws.Range("C23").Value = ' You\''ll finally learn what I\''m up\nto here. Meh-heh-heh.'
$ws.Range("D23").Value = ' Наконец-то вы узнаете зачем я\nздесь. Хе-хе-хе.'
$ws.Range("E23").Value = ' Îàëïîåø-óï âú ôèîàåóå èàœåí ÿ\nèäåòû. Öå-öå-öå.'

# --- Apply "group end" (thin-bottom-border) style to the last row of each dialogue block ---
$ws.Range("A6:E6").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)
$ws.Range("A6:E6").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)
$ws.Range("A6:E6").Copy()
$ws.Range("A23:E23").PasteSpecial(-4122)

# --- Row heights ---
$ws.Rows.Item(13).RowHeight = 21.6
$ws.Rows.Item(14).RowHeight = 21.6
$ws.Rows.Item(15).RowHeight = 31.8
$ws.Rows.Item(16).RowHeight = 21.6
$ws.Rows.Item(18).RowHeight = 21.6
$ws.Rows.Item(19).RowHeight = 21.6
$ws.Rows.Item(20).RowHeight = 21.6
$ws.Rows.Item(21).RowHeight = 31.8
$ws.Rows.Item(22).RowHeight = 31.8
$ws.Rows.Item(23).RowHeight = 21.6
$ws.Rows.Item(24).RowHeight = 42.0
$ws.Rows.Item(25).RowHeight = 42.0

# --- Update view / selection to match target ---
$ws.Application.Goto($ws.Range("A22"))
$ws.Range("B22").Select()

Write-Output "done"
